$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Cell = 'D2'; Value = '29.125.91'; ForceText = $false },
    @{ Cell = 'E2'; Value = '  +0.21%  '; ForceText = $false },
    @{ Cell = 'D3'; Value = '1.840.45'; ForceText = $false },
    @{ Cell = 'E3'; Value = '  +0.33%  '; ForceText = $false },
    @{ Cell = 'E4'; Value = '  +0.03%  '; ForceText = $false },
    @{ Cell = 'D5'; Value = '244.32'; ForceText = $true },
    @{ Cell = 'E5'; Value = '  -0.05%  '; ForceText = $false },
    @{ Cell = 'D6'; Value = '0.6255'; ForceText = $true },
    @{ Cell = 'E6'; Value = '  -1.47%  '; ForceText = $false },
    @{ Cell = 'E7'; Value = '  +0.11%  '; ForceText = $false },
    @{ Cell = 'D8'; Value = '0.07488'; ForceText = $true },
    @{ Cell = 'E8'; Value = '  -0.94%  '; ForceText = $false },
    @{ Cell = 'E9'; Value = '  -0.36%  '; ForceText = $false },
    @{ Cell = 'D10'; Value = '23.30'; ForceText = $true },
    @{ Cell = 'E10'; Value = '  +1.44%  '; ForceText = $false },
    @{ Cell = 'D11'; Value = '0.07714'; ForceText = $true },
    @{ Cell = 'E11'; Value = '  -0.35%  '; ForceText = $false },
    @{ Cell = 'D12'; Value = '1.880.82'; ForceText = $false },
    @{ Cell = 'E12'; Value = '  +1.52%  '; ForceText = $false },
    @{ Cell = 'D13'; Value = '5.017'; ForceText = $true },
    @{ Cell = 'E13'; Value = '  +0.06%  '; ForceText = $false },
    @{ Cell = 'D14'; Value = '0.6754'; ForceText = $true },
    @{ Cell = 'E14'; Value = '  +0.47%  '; ForceText = $false },
    @{ Cell = 'D15'; Value = '83.05'; ForceText = $true },
    @{ Cell = 'E15'; Value = '  -0.27%  '; ForceText = $false },
    @{ Cell = 'D16'; Value = '0.000009287'; ForceText = $true },
    @{ Cell = 'E16'; Value = '  -4.53%  '; ForceText = $false },
    @{ Cell = 'D17'; Value = '5.974'; ForceText = $true },
    @{ Cell = 'E17'; Value = '  -2.06%  '; ForceText = $false },
    @{ Cell = 'D18'; Value = '29.135.19'; ForceText = $false },
    @{ Cell = 'E18'; Value = '  +0.16%  '; ForceText = $false },
    @{ Cell = 'D19'; Value = '2.129.48'; ForceText = $false },
    @{ Cell = 'E19'; Value = '  +2.19%  '; ForceText = $false },
    @{ Cell = 'D20'; Value = '230.86'; ForceText = $true },
    @{ Cell = 'E20'; Value = '  +1.91%  '; ForceText = $false },
    @{ Cell = 'E21'; Value = '  +0.76%  '; ForceText = $false },
    @{ Cell = 'D22'; Value = '1.002'; ForceText = $true },
    @{ Cell = 'E22'; Value = '  +0.22%  '; ForceText = $false },
    @{ Cell = 'D23'; Value = '7.198'; ForceText = $true },
    @{ Cell = 'E23'; Value = '  -0.11%  '; ForceText = $false },
    @{ Cell = 'E24'; Value = '  +0.07%  '; ForceText = $false },
    @{ Cell = 'D25'; Value = '160.45'; ForceText = $true },
    @{ Cell = 'E25'; Value = '  -0.08%  '; ForceText = $false },
    @{ Cell = 'B26'; Value = 'Stellar'; ForceText = $false },
    @{ Cell = 'C26'; Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; ForceText = $false },
    @{ Cell = 'D26'; Value = '0.1394'; ForceText = $true },
    @{ Cell = 'E26'; Value = '  -0.78%  '; ForceText = $false },
    @{ Cell = 'B27'; Value = 'Cosmos'; ForceText = $false },
    @{ Cell = 'C27'; Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; ForceText = $false },
    @{ Cell = 'D27'; Value = '8.557'; ForceText = $true },
    @{ Cell = 'E27'; Value = '  -0.03%  '; ForceText = $false },
    @{ Cell = 'E28'; Value = '  -0.20%  '; ForceText = $false },
    @{ Cell = 'E29'; Value = '  +0.33%  '; ForceText = $false },
    @{ Cell = 'D30'; Value = '4.186'; ForceText = $true },
    @{ Cell = 'E30'; Value = '  +1.37%  '; ForceText = $false },
    @{ Cell = 'D31'; Value = '4.149'; ForceText = $true },
    @{ Cell = 'E31'; Value = '  +1.41%  '; ForceText = $false },
    @{ Cell = 'D32'; Value = '0.05584'; ForceText = $true },
    @{ Cell = 'E32'; Value = '  +3.41%  '; ForceText = $false },
    @{ Cell = 'D33'; Value = '1.207'; ForceText = $true },
    @{ Cell = 'E33'; Value = '  +0.23%  '; ForceText = $false },
    @{ Cell = 'D34'; Value = '0.7506'; ForceText = $true },
    @{ Cell = 'E34'; Value = '  +0.57%  '; ForceText = $false },
    @{ Cell = 'D35'; Value = '1.853'; ForceText = $true },
    @{ Cell = 'E35'; Value = '  -0.78%  '; ForceText = $false },
    @{ Cell = 'D36'; Value = '1.144'; ForceText = $true },
    @{ Cell = 'E36'; Value = '  -0.05%  '; ForceText = $false },
    @{ Cell = 'D37'; Value = '2.660'; ForceText = $true },
    @{ Cell = 'E37'; Value = '  +0.09%  '; ForceText = $false },
    @{ Cell = 'D38'; Value = '2.766'; ForceText = $true },
    @{ Cell = 'E38'; Value = '  +0.28%  '; ForceText = $false },
    @{ Cell = 'D39'; Value = '1.221.92'; ForceText = $false },
    @{ Cell = 'E39'; Value = '  -1.93%  '; ForceText = $false },
    @{ Cell = 'D41'; Value = '6.575'; ForceText = $true },
    @{ Cell = 'E41'; Value = '  -1.24%  '; ForceText = $false },
    @{ Cell = 'D42'; Value = '0.9039'; ForceText = $true },
    @{ Cell = 'E42'; Value = '  -0.40%  '; ForceText = $false },
    @{ Cell = 'D43'; Value = '1.001'; ForceText = $true },
    @{ Cell = 'E43'; Value = '  +0.11%  '; ForceText = $false },
    @{ Cell = 'D44'; Value = '2.022.35'; ForceText = $false },
    @{ Cell = 'E44'; Value = '  +1.88%  '; ForceText = $false },
    @{ Cell = 'D45'; Value = '102.23'; ForceText = $true },
    @{ Cell = 'E45'; Value = '  +0.24%  '; ForceText = $false },
    @{ Cell = 'D46'; Value = '66.24'; ForceText = $true },
    @{ Cell = 'E46'; Value = '  +1.91%  '; ForceText = $false },
    @{ Cell = 'D47'; Value = '0.00000000123'; ForceText = $true },
    @{ Cell = 'E47'; Value = '  +1.11%  '; ForceText = $false },
    @{ Cell = 'D48'; Value = '0.5093'; ForceText = $true },
    @{ Cell = 'E48'; Value = '  -0.40%  '; ForceText = $false },
    @{ Cell = 'D49'; Value = '0.4096'; ForceText = $true },
    @{ Cell = 'E49'; Value = '  -0.08%  '; ForceText = $false },
    @{ Cell = 'D50'; Value = '9.131'; ForceText = $true },
    @{ Cell = 'E50'; Value = '  +0.83%  '; ForceText = $false },
    @{ Cell = 'D51'; Value = '0.05844'; ForceText = $true },
    @{ Cell = 'E51'; Value = '  +1.23%  '; ForceText = $false }
)

foreach ($change in $changes) {
    $range = $ws.Range($change.Cell)
    if ($change.ForceText) {
        $range.NumberFormat = "@"
    }
    $range.Value = $change.Value
}
